$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F4").Value = 630
$ws1.Range("F5").Value = 175
$ws1.Range("F6").Value = 9481
$ws1.Range("F7").Value = 853
$ws1.Range("F8").Value = 333
$ws1.Range("F10").Value = 1178
$ws1.Range("F11").Value = 152
$ws1.Range("F12").Value = 102
$ws1.Range("F15").Value = 435
$ws1.Range("F16").Value = 94
$ws1.Range("F17").Value = 256
$ws1.Range("F18").Value = 1302

$ws4.Range("F4").Value = 14
$ws4.Range("F5").Value = 630
$ws4.Range("F6").Value = 175
$ws4.Range("F7").Value = 9481
$ws4.Range("F8").Value = 853
$ws4.Range("F9").Value = 333
$ws4.Range("F11").Value = 1178
$ws4.Range("F12").Value = 152
$ws4.Range("F13").Value = 102
$ws4.Range("F16").Value = 435
$ws4.Range("F17").Value = 94
$ws4.Range("F18").Value = 256
$ws4.Range("F19").Value = 1302
